$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "33.843.61"
Set-TextValue "E2" "  -0.99%  "

Set-TextValue "D3" "1.778.45"
Set-TextValue "E3" "  -1.46%  "

Set-TextValue "E4" "  +0.13%  "

Set-TextValue "D5" "223.96"
Set-TextValue "E5" "  +0.19%  "

Set-TextValue "D6" "0.546"
Set-TextValue "E6" "  -1.26%  "

Set-TextValue "E7" "  +0.13%  "

Set-TextValue "D8" "31.71"
Set-TextValue "E8" "  -2.14%  "

Set-TextValue "E9" "  +0.54%  "

Set-TextValue "D10" "0.0675"
Set-TextValue "E10" "  -6.37%  "

Set-TextValue "D11" "0.0935"
Set-TextValue "E11" "  +0.80%  "

Set-TextValue "D12" "2.036.55"
Set-TextValue "E12" "  -1.29%  "

Set-TextValue "D13" "11.15"
Set-TextValue "E13" "  +1.24%  "

Set-TextValue "D14" "1.767.15"
Set-TextValue "E14" "  -2.16%  "

Set-TextValue "D15" "33.865.81"
Set-TextValue "E15" "  -1.02%  "

Set-TextValue "D16" "0.608"
Set-TextValue "E16" "  -3.70%  "

Set-TextValue "E17" "  -2.41%  "

Set-TextValue "D18" "66.74"
Set-TextValue "E18" "  -2.78%  "

Set-TextValue "D19" "238.71"
Set-TextValue "E19" "  -3.78%  "

Set-TextValue "D20" "0.0₃0769"
Set-TextValue "E20" "  -2.85%  "

Set-TextValue "E21" "  +0.11%  "

Set-TextValue "D22" "10.54"
Set-TextValue "E22" "  -3.97%  "

Set-TextValue "E23" "  -2.74%  "

Set-TextValue "D24" "2.06"
Set-TextValue "E24" "  -3.11%  "

Set-TextValue "D25" "160.91"
Set-TextValue "E25" "  +0.76%  "

Set-TextValue "D26" "7.02"
Set-TextValue "E26" "  -1.40%  "

Set-TextValue "D27" "16.06"
Set-TextValue "E27" "  -3.58%  "

Set-TextValue "E28" "  -1.22%  "

Set-TextValue "E29" "  +0.27%  "

Set-TextValue "E30" "  +0.68%  "

Set-TextValue "E31" "  -3.31%  "

Set-TextValue "E32" "  -4.07%  "

Set-TextValue "E33" "  -0.16%  "

Set-TextValue "E34" "  -2.53%  "

Set-TextValue "D35" "1.389.60"
Set-TextValue "E35" "  -2.08%  "

Set-TextValue "D36" "0.632"
Set-TextValue "E36" "  -3.45%  "

Set-TextValue "E37" "  -1.66%  "

Set-TextValue "E38" "  -1.17%  "

Set-TextValue "E39" "  +4.18%  "

Set-TextValue "D40" "2.36"
Set-TextValue "E40" "  +0.00%  "

Set-TextValue "D41" "0.910"
Set-TextValue "E41" "  -3.80%  "

Set-TextValue "D42" "78.18"
Set-TextValue "E42" "  -3.08%  "

Set-TextValue "B43" "MXToken"
Set-TextValue "C43" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.65"
Set-TextValue "E43" "  -3.06%  "

Set-TextValue "B44" "InjectiveProtocol"
Set-TextValue "C44" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D44" "13.43"
Set-TextValue "E44" "  +11.65%  "

Set-TextValue "B45" "Kaspa"
Set-TextValue "C45" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D45" "0.0509"
Set-TextValue "E45" "  +2.39%  "

Set-TextValue "B46" "BabyDogeCoin"
Set-TextValue "C46" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D46" "0.0₆0138"
Set-TextValue "E46" "  +10.47%  "

Set-TextValue "E47" "  +2.59%  "

Set-TextValue "D48" "5.85"
Set-TextValue "E48" "  -1.72%  "

Set-TextValue "D49" "105.85"
Set-TextValue "E49" "  -2.16%  "

Set-TextValue "D50" "1.937.75"
Set-TextValue "E50" "  -1.45%  "

Set-TextValue "E51" "  +0.17%  "
